$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "building - refurbishment" category is being split into two new
# categories ("building - conventional renovation" and
# "building - serial renovation"). This adds one row to the id_action
# table, so first grow the table by one row via its ListObject so the
# table/autoFilter range resizes automatically.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Update B3 (was "building - refurbishment") to the first new category.
$ws.Range("B3").Value = "building - conventional renovation"

# Insert the second new category as a new row, and shift the technology
# rows down beneath it, renumbering id_action 2,3,4 -> 4,5,6.
$ws.Range("B4").Value = "building - serial renovation"
$ws.Range("B5").Value = "technology - new installation"
$ws.Range("B6").Value = "technology - similar change"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "technology - switch to new type of system"

# Update the view: zoom to 188% and move the selection to the new last row.
$excel.ActiveWindow.Zoom = 188
$ws.Range("B7").Select() | Out-Null
